$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Assigned" header in D1 (same style as the other headers)
$ws.Range("D1").Value = "Assigned"
$ws.Range("D1").Style = $ws.Range("C1").Style

# Assigned values for students in rows 2-12 (only first student assigned)
$ws.Range("D2").Value = $true
$ws.Range("D3").Value = $false
$ws.Range("D4").Value = $false
$ws.Range("D5").Value = $false
$ws.Range("D6").Value = $false
$ws.Range("D7").Value = $false
$ws.Range("D8").Value = $false
$ws.Range("D9").Value = $false
$ws.Range("D10").Value = $false
$ws.Range("D11").Value = $false
$ws.Range("D12").Value = $false

# Move selection to D3 (matches the author's cursor position after editing)
$ws.Range("D3").Select()
